$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing literal text storage (avoids
# Excel auto-converting numeric-looking strings like "1.001" or "30.616.85"
# into actual numbers), while preserving the cells original style.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

# Row 2
Set-TextValue "D2" "30.616.85"
$ws.Range("E2").Value = "  +0.73%  "

# Row 3
Set-TextValue "D3" "1.924.18"
$ws.Range("E3").Value = "  -0.10%  "

# Row 4
Set-TextValue "D4" "1.001"
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
Set-TextValue "D5" "246.86"
$ws.Range("E5").Value = "  +2.62%  "

# Row 6
Set-TextValue "D6" "1.001"
$ws.Range("E6").Value = "  +0.03%  "

# Row 7
Set-TextValue "D7" "0.4748"
$ws.Range("E7").Value = "  -0.33%  "

# Row 8
Set-TextValue "D8" "0.2888"
$ws.Range("E8").Value = "  +0.95%  "

# Row 9
Set-TextValue "D9" "0.06821"
$ws.Range("E9").Value = "  +3.50%  "

# Row 10
Set-TextValue "D10" "105.12"
$ws.Range("E10").Value = "  -1.02%  "

# Row 11
Set-TextValue "D11" "18.37"
$ws.Range("E11").Value = "  -3.87%  "

# Row 12
Set-TextValue "D12" "1.920.42"
$ws.Range("E12").Value = "  +0.22%  "

# Row 13
Set-TextValue "D13" "0.07681"
$ws.Range("E13").Value = "  +0.92%  "

# Row 14
Set-TextValue "D14" "5.337"
$ws.Range("E14").Value = "  +4.23%  "

# Row 15
$ws.Range("E15").Value = "  +1.56%  "

# Row 16
Set-TextValue "D16" "289.75"
$ws.Range("E16").Value = "  -4.03%  "

# Row 17
Set-TextValue "D17" "30.619.17"

# Row 18
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D18" "5.594"
$ws.Range("E18").Value = "  +6.35%  "

# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D19" "0.000007616"
$ws.Range("E19").Value = "  +1.65%  "

# Row 20
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D20" "1.001"
$ws.Range("E20").Value = "  +0.04%  "

# Row 21
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D21" "12.95"
$ws.Range("E21").Value = "  +0.13%  "

# Row 22
Set-TextValue "D22" "2.175.74"
$ws.Range("E22").Value = "  +0.38%  "

# Row 23
Set-TextValue "D23" "1.001"
$ws.Range("E23").Value = "  +0.10%  "

# Row 24
Set-TextValue "D24" "6.445"
$ws.Range("E24").Value = "  +2.17%  "

# Row 25
Set-TextValue "D25" "9.477"
$ws.Range("E25").Value = "  +2.95%  "

# Row 26
Set-TextValue "D26" "166.81"
$ws.Range("E26").Value = "  -1.05%  "

# Row 27
Set-TextValue "D27" "21.20"
$ws.Range("E27").Value = "  +7.23%  "

# Row 28
Set-TextValue "D28" "2.108"
$ws.Range("E28").Value = "  +5.38%  "

# Row 29
$ws.Range("E29").Value = "  -4.34%  "

# Row 30
Set-TextValue "D30" "1.401"
$ws.Range("E30").Value = "  +3.64%  "

# Row 31
Set-TextValue "D31" "4.184"
$ws.Range("E31").Value = "  +2.39%  "

# Row 32
Set-TextValue "D32" "4.045"
$ws.Range("E32").Value = "  +3.22%  "

# Row 33
Set-TextValue "D33" "0.05026"
$ws.Range("E33").Value = "  +0.53%  "

# Row 34
Set-TextValue "D34" "0.7326"
$ws.Range("E34").Value = "  -0.99%  "

# Row 35
$ws.Range("E35").Value = "  -0.40%  "

# Row 36
Set-TextValue "D36" "0.02058"
$ws.Range("E36").Value = "  +5.74%  "

# Row 37
Set-TextValue "D37" "2.737"
$ws.Range("E37").Value = "  +0.24%  "

# Row 38
Set-TextValue "D38" "2.689"
$ws.Range("E38").Value = "  -0.38%  "

# Row 39
Set-TextValue "D39" "2.048"
$ws.Range("E39").Value = "  -0.26%  "

# Row 40
Set-TextValue "D40" "111.58"
$ws.Range("E40").Value = "  +3.98%  "

# Row 41
Set-TextValue "D41" "0.8724"
$ws.Range("E41").Value = "  -0.04%  "

# Row 42
Set-TextValue "D42" "0.4380"
$ws.Range("E42").Value = "  +6.25%  "

# Row 43
Set-TextValue "D43" "5.917"
$ws.Range("E43").Value = "  +2.32%  "

# Row 45
Set-TextValue "D45" "67.71"
$ws.Range("E45").Value = "  -2.80%  "

# Row 46
Set-TextValue "D46" "7.283"

# Row 47
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue "D47" "48.63"
$ws.Range("E47").Value = "  +16.45%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "9.307"
$ws.Range("E48").Value = "  +0.55%  "

# Row 49
Set-TextValue "D49" "0.1243"
$ws.Range("E49").Value = "  +3.30%  "

# Row 50
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue "D50" "34.93"
$ws.Range("E50").Value = "  +0.39%  "

# Row 51
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D51" "1.463"
$ws.Range("E51").Value = "  +7.03%  "
